$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows of this weekly price table were re-shuffled: for each row, the
# "Fecha" date together with its associated Volumen/Precio/Unidad/Origen data
# (columns D, M, N, O, P, Q, R, S, T) now belongs to a different row. The
# other descriptive columns (A,B,C,E,F,G,H,I,J,K,L) are identical for every
# row and stay untouched.
$cols = @("D","M","N","O","P","Q","R","S","T")

# Snapshot the moving columns for rows 2-22 before any writes, since the
# permutation reassigns each row from a different source row and later
# writes must not read already-overwritten data.
$snapshot = @{}
for ($r = 2; $r -le 22; $r++) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: target row number -> source row number (the row whose data should
# end up in the target row).
$mapping = @{
    2 = 17
    3 = 3
    4 = 5
    5 = 22
    6 = 11
    7 = 10
    8 = 14
    9 = 2
    10 = 19
    11 = 16
    12 = 9
    13 = 12
    14 = 15
    15 = 8
    16 = 7
    17 = 18
    18 = 13
    19 = 21
    20 = 20
    21 = 6
    22 = 4
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $src = $snapshot[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value = $src[$c]
    }
}

Write-Host "Row data reshuffled"
